# Swate RNASeq computational-analysis template:
# Remove the ER (endpoint repository) tag values from the
# "SwateTemplateMetadata" sheet — this template does not target an ER,
# so the "ER" / "ER Term Accession Number" / "ER Term Source REF" row
# values ("GEO" / the DPBO accession URL / "DPBO") are cleared out,
# matching the commit message "removed ER tags from non-ER templates
# and non-ER tags".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

$ws.Range("B8:B10").ClearContents() | Out-Null

# Leave the cursor where the author last left it while editing this sheet.
$ws.Range("G10").Select() | Out-Null
